$wb = $excel.ActiveWorkbook

# Sheet 1 = 檢核表 (unchanged), Sheet 2 = 未銷帳 (existing, to be duplicated)
$wsCheck  = $wb.Worksheets.Item(1)
$wsUnsold = $wb.Worksheets.Item(2)

# Duplicate the "未銷帳" sheet and place the copy right after "檢核表"
# (Copy's 2nd arg = "After" target sheet)
$wsUnsold.Copy([System.Reflection.Missing]::Value, $wsCheck) | Out-Null

# The freshly inserted copy is now in position 2; rename it
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "應收應付"
$wsNew.Activate() | Out-Null
$wsNew.Range("A2").Select() | Out-Null

# Re-resolve the original "未銷帳" sheet by name (its index shifted to 3)
# and move its selection to A2 as well
$wsOrig = $wb.Worksheets.Item("未銷帳")
$wsOrig.Activate() | Out-Null
$wsOrig.Range("A2").Select() | Out-Null

# Restore "檢核表" as the active/selected tab (it was untouched in the source edit)
$wsCheck.Activate() | Out-Null
